$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Reorder sheets: move "SV TEF" to sit before "SVP" (so order becomes
#    ... SV PAF, SV TEF, SVP). This also swaps which physical sheet part
#    backs "SVP" vs "SV TEF" internally, matching the target workbook.xml
#    (SV TEF keeps sheetId 6 / gets r:id rId5, SVP keeps sheetId 5 / gets
#    r:id rId6).
# ---------------------------------------------------------------------------
$wsTEF = $wb.Worksheets.Item("SV TEF")
$wsSVP = $wb.Worksheets.Item("SVP")
$wsTEF.Move($wsSVP)

# ---------------------------------------------------------------------------
# 2. Append the new Ultra Rare / SVP promo rows to the "SVP" sheet.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SVP")

$newRows = @(
    @(86,  "Mabosstiff ex", "Mabosstiff",   "Darkness",  "SV_Promo", "ex_SV"),
    @(103, "Houndoom ex",   "Houndoom",     "Darkness",  "SV_Promo", "ex_SV"),
    @(104, "Melmetal ex",   "Melmetal",     "Metal",     "SV_Promo", "ex_SV"),
    @(87,  "Sprigatito ex", "Sprigatito",   "Grass",     "SV_Promo", "ex_SV"),
    @(88,  "Pikachu ex",    "Pikachu",      "Lightning", "SV_Promo", "Full_Art_Pokemon_SV"),
    @(97,  "Flutter Mane",  "Flutter_Mane", "Psychic",   "SV_Promo", "Special_Art_Pokemon_SV_Ancient"),
    @(98,  "Iron Thorns",   "Iron_Thorns",  "Lightning", "SV_Promo", "Special_Art_Pokemon_SV_Future")
)

$startRow = 17
$row = $startRow
foreach ($data in $newRows) {
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]
    $ws.Cells.Item($row, 5).Value = $data[4]
    $ws.Cells.Item($row, 6).Value = $data[5]
    $row = $row + 1
}
$endRow = $row - 1

# Fill column G with the same formula as the existing rows, letting Excel
# adjust the row-relative references for each destination row. The formula
# text must reference the first row of the destination range (startRow) so
# the COM layer anchors the relative references correctly.
$formula = '="new Card(""" & B' + $startRow + ' & """, Pokedex." & C' + $startRow + ' & ", Rarity." & F' + $startRow + ' & ", Types." & D' + $startRow + ' & ", Sets." & E' + $startRow + ' & ", " & A' + $startRow + ' & "),"'
$ws.Range("G" + $startRow + ":G" + $endRow).Formula = $formula

# ---------------------------------------------------------------------------
# 3. Update the view/selection so SVP is the active tab with G22:G23 selected.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("G22:G23").Select()
